# Refine API specs and other misc changes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: /api/rooms -> add /api/rooms/:id row, fix D10 style ---
$ws.Range("A10").Value = "/api/rooms/:id"
$ws.Range("C10").Value = "{data:{room:{}}}"
# D10 style normalized to match the rest of the row (s=0 -> s=2)
$ws.Range("E10").Copy()
$ws.Range("D10").PasteSpecial(-4122)

# --- Row 11: new /api/rooms/create details row ---
$ws.Range("A11").Value = "/api/rooms/create"
$ws.Range("B11").Value = "{number,baseRent}"
$ws.Range("C11").Value = "{data:{room:{}}}"
$ws.Range("D11").Value = "The newly created Room"

# --- Row 12: clear out the old (duplicated) /api/tenants content ---
$ws.Range("A12").ClearContents()
$ws.Range("B12").ClearContents()
$ws.Range("C12").ClearContents()
$ws.Range("D12").ClearContents()

# --- Row 13: /api/tenants ---
$ws.Range("A13").Value = "/api/tenants"
$ws.Range("C13").Value = "{data:{tenants:[]}}"

# --- Row 14: /api/tenants/:id ---
$ws.Range("A14").Value = "/api/tenants/:id"
$ws.Range("C14").Value = "{data:{tenant:{}}}"

# --- Row 15: /api/tenants/create details ---
$ws.Range("A15").Value = "/api/tenants/create"
$ws.Range("B15").Value = "{name,phoneNumber,aadharCard,room}"
$ws.Range("C15").Value = "{data:{tenant:{}}}"
$ws.Range("D15").Value = "The newly created Tenant"

# --- Row 17: /api/transactions/?room, fix D17 style ---
$ws.Range("A17").Value = "/api/transactions/?room"
$ws.Range("C17").Value = "{data:{transactions:[]}}"
# D17 style normalized to match the rest of the row (s=0 -> s=2)
$ws.Range("E17").Copy()
$ws.Range("D17").PasteSpecial(-4122)

# --- Row 18: /api/transactions/:id ---
$ws.Range("A18").Value = "/api/transactions/:id"
$ws.Range("C18").Value = "{data:{transaction:{}}}"

# --- Row 19: transaction create details ---
$ws.Range("A19").Value = "/api/tenants/create/?roomNumber"
$ws.Range("B19").Value = "{room,balance,transfer,remarks}"
$ws.Range("C19").Value = "{data:{transaction:{}}}"
$ws.Range("D19").Value = "The newly created Transaction"

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 36.75
$ws.Columns.Item(2).ColumnWidth = 40.75

# --- Extend used range with two new blank rows (33, 34), matching row 32's formatting ---
$ws.Range("A32:AI32").Copy()
$ws.Range("A33:AI33").PasteSpecial(-4122)
$ws.Range("A32:AI32").Copy()
$ws.Range("A34:AI34").PasteSpecial(-4122)

# --- Update the active selection ---
$null = $ws.Range("B12").Select()
